# Product Backlog update: "cap nhat backlog + code"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")
$ws.Activate()

# --- User Login (row 9) ---
$ws.Range("C9").Value = "Medium"
$ws.Range("D9").Value = 220
$ws.Range("I9").Value = "Done"
$ws.Range("K9").Value = "Done"
$ws.Range("L9").Value = "Done"
$ws.Range("M9").Value = 90

# --- User Register (row 10) ---
$ws.Range("C10").Value = "Medium"
$ws.Range("D10").Value = 220
$ws.Range("I10").Value = "Done"
$ws.Range("K10").Value = "Done"
$ws.Range("L10").Value = "Done"
$ws.Range("M10").Value = 90

# --- Home Page (row 11) ---
$ws.Range("C11").Value = "Medium"
$ws.Range("D11").Value = 220
$ws.Range("I11").Value = "Done"
$ws.Range("M11").Value = 90

# --- Products List (Customer) (row 12) ---
$ws.Range("I12").Value = "Done"
$ws.Range("M12").Value = 70

# --- Cart Details (row 13) ---
$ws.Range("I13").Value = "Done"

# --- Admin Dashboard (row 14) ---
$ws.Range("I14").Value = "Done"
$ws.Range("M14").Value = 70

# --- Users List (row 15) ---
$ws.Range("D15").Value = 240
$ws.Range("I15").Value = "Done"
$ws.Range("M15").Value = 70

# --- Order Information (overview) (row 16) ---
$ws.Range("D16").Value = 170
$ws.Range("I16").Value = "Done"

# --- Rows 17-33: SRS Status column moves from Pending to Done ---
$ws.Range("I17:I33").Value = "Done"

# --- Compare Product (row 34): SRS Status moves from Pending to Doing ---
$ws.Range("I34").Value = "Doing"

# --- Restore view/selection state to match where the edits were made ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I30").Select()
